$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 36.94436433333333
$ws.Range("H2").Value = 110.833093
$ws.Range("I2").Value = 0.8328964975864823
$ws.Range("J2").Value = 0.8328964975864824
$ws.Range("M2").Value = 0.415892
$ws.Range("N2").Value = 1.247676
$ws.Range("Q2").Value = 15.36486557131867
$ws.Range("R2").Value = 138.283790141868
$ws.Range("S2").Value = 0.8328964975864823
$ws.Range("T2").Value = 0.8328964975864824

# Row 3
$ws.Range("I3").Value = 0.07608399754092349
$ws.Range("J3").Value = 0.07608399754092349
$ws.Range("M3").Value = 0.415892
$ws.Range("N3").Value = 1.247676
$ws.Range("Q3").Value = 1.403560223548
$ws.Range("R3").Value = 12.632042011932
$ws.Range("S3").Value = 0.07608399754092349
$ws.Range("T3").Value = 0.07608399754092349

# Row 4
$ws.Range("G4").Value = 4.037305666666668
$ws.Range("H4").Value = 12.111917
$ws.Range("I4").Value = 0.09101950487259411
$ws.Range("J4").Value = 0.09101950487259411
$ws.Range("M4").Value = 0.415892
$ws.Range("N4").Value = 1.247676
$ws.Range("Q4").Value = 1.679083128321334
$ws.Range("R4").Value = 15.111748154892
$ws.Range("S4").Value = 0.09101950487259411
$ws.Range("T4").Value = 0.09101950487259411
